$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.020.98'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '1.827.73'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9989'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '240.50'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.54%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.6193'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -6.73%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '44.45'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +6.06%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07483'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.99%  '
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '22.64'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.06%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07614'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.81%  '
$ws.Range('D13').Value = '1.825.86'
$ws.Range('E13').Value = '  -0.91%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.943'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.86%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6617'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.93%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '81.88'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.17%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000009043'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +8.19%  '
$ws.Range('E18').Value = '  -2.24%  '
$ws.Range('D19').Value = '29.015.39'
$ws.Range('E19').Value = '  -0.52%  '
$ws.Range('D20').Value = '2.076.54'
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '223.92'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.82%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '12.31'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.24%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '7.160'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.24%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.000'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '159.47'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.392'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.54%  '
$ws.Range('E28').Value = '  -4.42%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '17.77'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.14%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.492'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.29%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.025'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.204'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.04%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.034'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.83%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05204'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.87%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.824'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.21%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.148'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.75%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.7309'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.13%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.644'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('D39').Value = '1.273.99'
$ws.Range('E39').Value = '  +0.28%  '
$ws.Range('E40').Value = '  +0.60%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.01779'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.05%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.311'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +7.25%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.8907'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -4.37%  '
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('D46').Value = '1.974.88'
$ws.Range('E46').Value = '  -0.68%  '
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('E49').Value = '  -0.81%  '
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.3950'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.67%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.677'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -4.68%  '
